# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with the latest scraped values (GitHub Actions data-refresh job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.584.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4702"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2750"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06352"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.841.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.158"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6316"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.583.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "243.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007374"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.90%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.325"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.891"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1017"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.12%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.058"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.870"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04920"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.150"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7083"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.709"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01913"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8799"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.992"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("E41").Value = "  +0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4090"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.551"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.278"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.60%  "

$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.624"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.46%  "

$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05543"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3706"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.51%  "
